$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Lao" column (D) and "Total" column (G) values for rows 2 (Meteors) and 3 (Successes)
$ws.Range("D2").Value = 1758
$ws.Range("G2").Value = 1758
$ws.Range("D3").Value = 13
$ws.Range("G3").Value = 13
